# Populate the ASM / Archive disk-usage percentage readings for
# 10/13 (P), 10/14 (Q) and 10/15 (R) across the monitored rows.
#
# These cells must end up as plain text (shared-string) values like the
# existing "43%" / "O92" style readings in column O, NOT as numeric
# percentages - so we can't just do $range.Value = "42%" (Excel would
# auto-convert that into 0.42 formatted as a percent). Instead we copy
# the text out of a cell that already holds the exact string we need
# (Copy preserves the text data type), and for values that do not yet
# exist anywhere in the sheet we stage them in a scratch cell that is
# already Text-formatted, copy that into place, then restore the
# destination's original look (fill/border) via PasteSpecial of formats
# from a same-column neighbour, and finally clean the scratch cell back
# to empty.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$scratch = $ws.Range("D12")

function Set-TextFromSource([string]$destAddr, [string]$srcAddr, [string]$styleRefAddr) {
    $dst = $ws.Range($destAddr)
    $src = $ws.Range($srcAddr)
    $src.Copy($dst)

    if ($styleRefAddr) {
        $styleRef = $ws.Range($styleRefAddr)
        $styleRef.Copy()
        $dst.PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
    }
}

function Set-NewTextValue([string]$destAddr, [string]$text, [string]$styleRefAddr) {
    $dst = $ws.Range($destAddr)
    $styleRef = $ws.Range($styleRefAddr)

    $scratch.Value = $text
    $scratch.Copy($dst)

    $styleRef.Copy()
    $dst.PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

    $scratch.ClearContents()
}

# --- Row 8 : DBERP1 - Archive Volume(used(%) ---
Set-NewTextValue  "P8" "42%" "O8"
Set-NewTextValue  "Q8" "47%" "O8"
Set-TextFromSource "R8" "G8" "O8"

# --- Row 19 : DBPOS1 - Archive Volume(used(%) ---
Set-NewTextValue  "P19" "50%" "O19"
Set-TextFromSource "Q19" "O19"
Set-TextFromSource "R19" "M19"

# --- Row 30 : GPOS1 - Archive Volume(used(%) ---
Set-NewTextValue  "P30" "28%" "O30"
Set-TextFromSource "Q30" "P30"
Set-TextFromSource "R30" "P30"

# --- Row 40 : OGG - Oradata01 ---
Set-TextFromSource "P40" "O40"
Set-TextFromSource "Q40" "O40"
Set-TextFromSource "R40" "O40"

# --- Row 41 : OGG - Oradata02 ---
Set-TextFromSource "P41" "O41"
Set-TextFromSource "Q41" "O41"
Set-TextFromSource "R41" "O41"

# --- Row 42 : OGG - Oradata03 ---
Set-TextFromSource "P42" "O42"
Set-TextFromSource "Q42" "O42"
Set-TextFromSource "R42" "O42"

# --- Row 43 : OGG - Oradata04 ---
Set-TextFromSource "P43" "O43"
Set-TextFromSource "Q43" "O43"
Set-TextFromSource "R43" "O43"

# --- Row 44 : OGG - Oradata05 ---
Set-TextFromSource "P44" "E44" "O44"
Set-TextFromSource "Q44" "O44"
Set-TextFromSource "R44" "O44"

# --- Row 45 : OGG - Oradata06 ---
Set-TextFromSource "P45" "O45"
Set-TextFromSource "Q45" "O45"
Set-TextFromSource "R45" "O45"
